$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.956.88"
$ws.Range("E2").Value = "  -0.37%  "
$ws.Range("D3").Value = "1.675.82"
$ws.Range("E3").Value = "  +0.02%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.07"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.55%  "
$ws.Range("E6").Value = "  +1.43%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.252"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.45%  "
$ws.Range("E9").Value = "  +0.22%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.33"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.59%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0887"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.43%  "
$ws.Range("D12").Value = "1.911.87"
$ws.Range("E12").Value = "  +0.00%  "
$ws.Range("D13").Value = "1.694.32"
$ws.Range("E13").Value = "  +1.10%  "
$ws.Range("E14").Value = "  -0.26%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.528"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.93%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.80"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.19%  "
$ws.Range("D17").Value = "26.960.37"
$ws.Range("E17").Value = "  -0.39%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "236.93"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.35%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.03"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.64%  "
$ws.Range("E20").Value = "  -0.97%  "
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("E22").Value = "  -0.85%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.20"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.11%  "
$ws.Range("E24").Value = "  -2.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.75"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.26"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.62%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.05"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.36%  "
$ws.Range("E28").Value = "  -1.35%  "
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0497"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.33%  "
$ws.Range("E31").Value = "  -0.64%  "
$ws.Range("E32").Value = "  -0.09%  "
$ws.Range("D33").Value = "1.485.70"
$ws.Range("E33").Value = "  +0.74%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.15"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.92%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.67"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.54%  "
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.587"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0173"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.70%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.896"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.56%  "
$ws.Range("E40").Value = "  -3.31%  "
$ws.Range("E41").Value = "  +0.90%  "
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("E43").Value = "  +1.98%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "67.29"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.02%  "
$ws.Range("D45").Value = "1.816.77"
$ws.Range("E45").Value = "  -0.45%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.779"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.49"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.13%  "
$ws.Range("E48").Value = "  -0.53%  "
$ws.Range("E49").Value = "  +1.67%  "
$ws.Range("E50").Value = "  +0.65%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.69"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.62%  "
